$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Cells.Item(1, 26)  # Z1: unused scratch cell for text-safe writes

$scratch.Formula = '="66.700.97"'
$scratch.Copy()
$ws.Cells.Item(2, 4).PasteSpecial(-4163)
$ws.Cells.Item(2, 5).Value = '  -3.96%  '
$scratch.Formula = '="3.336.22"'
$scratch.Copy()
$ws.Cells.Item(3, 4).PasteSpecial(-4163)
$ws.Cells.Item(3, 5).Value = '  -1.35%  '
$ws.Cells.Item(4, 5).Value = '  -0.04%  '
$scratch.Formula = '="572.90"'
$scratch.Copy()
$ws.Cells.Item(5, 4).PasteSpecial(-4163)
$ws.Cells.Item(5, 5).Value = '  -3.36%  '
$scratch.Formula = '="181.68"'
$scratch.Copy()
$ws.Cells.Item(6, 4).PasteSpecial(-4163)
$ws.Cells.Item(6, 5).Value = '  -5.14%  '
$ws.Cells.Item(7, 5).Value = '  -0.05%  '
$ws.Cells.Item(8, 5).Value = '  -1.03%  '
$scratch.Formula = '="0.129"'
$scratch.Copy()
$ws.Cells.Item(9, 4).PasteSpecial(-4163)
$ws.Cells.Item(9, 5).Value = '  -3.84%  '
$ws.Cells.Item(10, 5).Value = '  -1.73%  '
$scratch.Formula = '="0.401"'
$scratch.Copy()
$ws.Cells.Item(11, 4).PasteSpecial(-4163)
$ws.Cells.Item(11, 5).Value = '  -4.49%  '
$scratch.Formula = '="3.916.33"'
$scratch.Copy()
$ws.Cells.Item(12, 4).PasteSpecial(-4163)
$ws.Cells.Item(12, 5).Value = '  -1.45%  '
$ws.Cells.Item(13, 5).Value = '  -1.78%  '
$scratch.Formula = '="27.06"'
$scratch.Copy()
$ws.Cells.Item(14, 4).PasteSpecial(-4163)
$ws.Cells.Item(14, 5).Value = '  -5.66%  '
$scratch.Formula = '="66.773.68"'
$scratch.Copy()
$ws.Cells.Item(15, 4).PasteSpecial(-4163)
$ws.Cells.Item(15, 5).Value = '  -3.95%  '
$scratch.Formula = '="0.0000167"'
$scratch.Copy()
$ws.Cells.Item(16, 4).PasteSpecial(-4163)
$ws.Cells.Item(16, 5).Value = '  -2.71%  '
$scratch.Formula = '="3.352.18"'
$scratch.Copy()
$ws.Cells.Item(17, 4).PasteSpecial(-4163)
$ws.Cells.Item(17, 5).Value = '  -1.05%  '
$scratch.Formula = '="437.36"'
$scratch.Copy()
$ws.Cells.Item(18, 4).PasteSpecial(-4163)
$ws.Cells.Item(18, 5).Value = '  -2.62%  '
$ws.Cells.Item(19, 2).Value = 'Polkadot'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$scratch.Formula = '="5.68"'
$scratch.Copy()
$ws.Cells.Item(19, 4).PasteSpecial(-4163)
$ws.Cells.Item(19, 5).Value = '  -2.93%  '
$ws.Cells.Item(20, 2).Value = 'Chainlink'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$scratch.Formula = '="13.57"'
$scratch.Copy()
$ws.Cells.Item(20, 4).PasteSpecial(-4163)
$ws.Cells.Item(20, 5).Value = '  -1.79%  '
$scratch.Formula = '="7.59"'
$scratch.Copy()
$ws.Cells.Item(21, 4).PasteSpecial(-4163)
$ws.Cells.Item(21, 5).Value = '  -2.92%  '
$scratch.Formula = '="73.69"'
$scratch.Copy()
$ws.Cells.Item(22, 4).PasteSpecial(-4163)
$ws.Cells.Item(22, 5).Value = '  -1.91%  '
$ws.Cells.Item(23, 5).Value = '  -0.14%  '
$ws.Cells.Item(24, 5).Value = '  -0.82%  '
$ws.Cells.Item(25, 5).Value = '  -4.15%  '
$ws.Cells.Item(26, 5).Value = '  +0.03%  '
$ws.Cells.Item(27, 5).Value = '  -4.91%  '
$ws.Cells.Item(28, 5).Value = '  -0.26%  '
$ws.Cells.Item(29, 5).Value = '  -1.62%  '
$scratch.Formula = '="22.85"'
$scratch.Copy()
$ws.Cells.Item(30, 4).PasteSpecial(-4163)
$ws.Cells.Item(30, 5).Value = '  -2.53%  '
$scratch.Formula = '="5.29"'
$scratch.Copy()
$ws.Cells.Item(31, 4).PasteSpecial(-4163)
$ws.Cells.Item(31, 5).Value = '  -6.50%  '
$ws.Cells.Item(32, 5).Value = '  -0.04%  '
$scratch.Formula = '="6.80"'
$scratch.Copy()
$ws.Cells.Item(33, 4).PasteSpecial(-4163)
$ws.Cells.Item(33, 5).Value = '  -2.82%  '
$ws.Cells.Item(34, 5).Value = '  -4.54%  '
$scratch.Formula = '="162.14"'
$scratch.Copy()
$ws.Cells.Item(35, 4).PasteSpecial(-4163)
$ws.Cells.Item(35, 5).Value = '  -1.95%  '
$ws.Cells.Item(36, 5).Value = '  -4.35%  '
$scratch.Formula = '="27.84"'
$scratch.Copy()
$ws.Cells.Item(37, 4).PasteSpecial(-4163)
$ws.Cells.Item(37, 5).Value = '  +1.25%  '
$ws.Cells.Item(38, 5).Value = '  -5.68%  '
$scratch.Formula = '="2.820.91"'
$scratch.Copy()
$ws.Cells.Item(39, 4).PasteSpecial(-4163)
$ws.Cells.Item(39, 5).Value = '  +2.33%  '
$ws.Cells.Item(40, 5).Value = '  -3.15%  '
$scratch.Formula = '="4.42"'
$scratch.Copy()
$ws.Cells.Item(41, 4).PasteSpecial(-4163)
$ws.Cells.Item(41, 5).Value = '  -3.70%  '
$scratch.Formula = '="6.24"'
$scratch.Copy()
$ws.Cells.Item(42, 4).PasteSpecial(-4163)
$ws.Cells.Item(42, 5).Value = '  -5.13%  '
$scratch.Formula = '="0.0668"'
$scratch.Copy()
$ws.Cells.Item(44, 4).PasteSpecial(-4163)
$ws.Cells.Item(44, 5).Value = '  -3.42%  '
$ws.Cells.Item(45, 5).Value = '  -4.04%  '
$scratch.Formula = '="2.35"'
$scratch.Copy()
$ws.Cells.Item(46, 4).PasteSpecial(-4163)
$ws.Cells.Item(46, 5).Value = '  -7.59%  '
$scratch.Formula = '="321.02"'
$scratch.Copy()
$ws.Cells.Item(47, 4).PasteSpecial(-4163)
$ws.Cells.Item(47, 5).Value = '  -5.85%  '
$ws.Cells.Item(48, 5).Value = '  -4.12%  '
$scratch.Formula = '="0.987"'
$scratch.Copy()
$ws.Cells.Item(49, 4).PasteSpecial(-4163)
$ws.Cells.Item(49, 5).Value = '  -3.84%  '
$scratch.Formula = '="6.16"'
$scratch.Copy()
$ws.Cells.Item(50, 4).PasteSpecial(-4163)
$ws.Cells.Item(50, 5).Value = '  -2.91%  '
$ws.Cells.Item(51, 5).Value = '  -1.57%  '

$scratch.Clear()
$excel.CutCopyMode = $false
